$wb = $excel.ActiveWorkbook

# --- Update the text note on "Hoja1" sheet (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.85 = 6700.91 pesos`n✅ 6700.91 pesos = 1.84 = 955.27 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 541.995
$wsTasas.Range("O10").Value = 3631.86
$wsTasas.Range("N12").Value = 3648
$wsTasas.Range("O12").Value = 520.05
